$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the win/loss/tie record columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from the existing
# last header cell (AC1) onto the new header cells so they match style "1".
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team win/loss/tie record for every player row (2 through 48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 65   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 97   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
